# Fruta / hortaliza, semanal
# Insert a new weekly data row at the top of the data block (row 2),
# pushing all existing data rows down by one, then populate the new
# row with this week's reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift all existing data rows (2..61) down by one (to 3..62),
# carrying values + formatting (keeps the date-formatted D column intact).
$src = $ws.Range("A2:R61")
$dst = $ws.Range("A3:R62")
$src.Copy($dst)

# Populate the newly opened row 2 with this week's record.
$ws.Range("D2").Value2  = 44756
$ws.Range("J2").Value2  = 240
$ws.Range("K2").Value2  = 30000
$ws.Range("L2").Value2  = 32000
$ws.Range("M2").Value2  = 31000
$ws.Range("P2").Value2  = 1240

Write-Output "done"
